$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.65
$ws.Range("A6").Value = -22.056
$ws.Range("A7").Value = -21.304
$ws.Range("B7").Value = 5.872
$ws.Range("B12").Value = 4.786
$ws.Range("D12").Value = -7.392999999999999
$ws.Range("C13").Value = -13.322
$ws.Range("C14").Value = -11.992
$ws.Range("B15").Value = 5.087000000000001
$ws.Range("A16").Value = -21.771
$ws.Range("C16").Value = -13.308
$ws.Range("C19").Value = -12.228
$ws.Range("A20").Value = -21.17
$ws.Range("B20").Value = 6.672
$ws.Range("B21").Value = 8.904
$ws.Range("B22").Value = 7.485000000000001
$ws.Range("C22").Value = -13.29
$ws.Range("D22").Value = -8.063000000000001
$ws.Range("B23").Value = 7.359999999999999
$ws.Range("A28").Value = -21.687
$ws.Range("A29").Value = -21.648
$ws.Range("B29").Value = 5.625999999999999
$ws.Range("D29").Value = -7.170999999999999
$ws.Range("A32").Value = -21.486
$ws.Range("B34").Value = 7.821
$ws.Range("D34").Value = -7.997999999999999
$ws.Range("C36").Value = -12.683
$ws.Range("A40").Value = -20.25
$ws.Range("B42").Value = 7.583000000000001
$ws.Range("B43").Value = 5.907000000000001
$ws.Range("D43").Value = -8.145000000000001
$ws.Range("B44").Value = 5.171
$ws.Range("B45").Value = 5.048
$ws.Range("A46").Value = -20.773
$ws.Range("B46").Value = 7.259
$ws.Range("C46").Value = -13.976
$ws.Range("D48").Value = -7.797
$ws.Range("B50").Value = 5.291
$ws.Range("C50").Value = -13.132
$ws.Range("A51").Value = -20.891
$ws.Range("B51").Value = 6.7
$ws.Range("A52").Value = -21.561
$ws.Range("A57").Value = -22.074
$ws.Range("A59").Value = -22.248
$ws.Range("D60").Value = -8.245000000000001
$ws.Range("A62").Value = -21.956
$ws.Range("A66").Value = -21.54
$ws.Range("B66").Value = 6.422
$ws.Range("B67").Value = 5.249
$ws.Range("D68").Value = -6.901000000000001
$ws.Range("D70").Value = -7.103
$ws.Range("A73").Value = -19.95800000000001
$ws.Range("D73").Value = -8.363000000000001
$ws.Range("A74").Value = -21.061
$ws.Range("B79").Value = 5.529999999999999
$ws.Range("B84").Value = 5.674
$ws.Range("D87").Value = -8.192
$ws.Range("A92").Value = -21.333
$ws.Range("B92").Value = 5.499000000000001
$ws.Range("D92").Value = -6.278999999999999
$ws.Range("C95").Value = -11.719
$ws.Range("B97").Value = 7.025999999999999
$ws.Range("C97").Value = -13.631
$ws.Range("A100").Value = -21.678
$ws.Range("D101").Value = -8.188999999999998
